$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Update "Contenu du stage" counts: one fewer C# student, one more ASSEMBLEUR student
$ws.Range("E16").Value = 19
$ws.Range("E19").Value = 3

# Update the matching displayed percentages (plain text cells, not formulas).
# Force text so Excel's auto percent-number heuristic doesn't convert "76 %" / "12 %"
# into a numeric 0.76 / 0.12, then clear the format back to the default (unstyled)
# cell so no stray number-format style is introduced.
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "76 %"
$ws.Range("G16").ClearFormats()

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "12 %"
$ws.Range("G19").ClearFormats()
